$wb = $excel.ActiveWorkbook

# --- Rename sheets (Step1 -> Employee sheet, Step1_Des -> Description) ---
$wb.Worksheets.Item("Step1").Name = "Employee sheet"
$wb.Worksheets.Item("Step1_Des").Name = "Description"

# --- Update the "Objectives:" header on the Description sheet to "Description:" ---
$wsDesc = $wb.Worksheets.Item("Description")
$wsDesc.Activate()
$wsDesc.Range("C4").Value = "Description:"

# Restore the view/selection on the Description sheet to match the saved state
$wsDesc.Range("C22:C26").Select()

# --- Restore the view/selection on the Step2 sheet (left active, matches activeTab) ---
$wsStep2 = $wb.Worksheets.Item("Step2")
$wsStep2.Activate()
$wsStep2.Range("A136:XFD1048576,A1:XFD135").Select()
